$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data region (rows 3-6) which had the original temperature/capacity/color/comment table
$ws.Range("A3:D6").ClearContents()

# New data: sample1..sample4 with numeric + comment columns
$data = @(
    @("sample1", 25.5, 25, "comment 1"),
    @("sample2", 23.5, 23, "comment 2"),
    @("sample3", 10.5, 10, "comment 3"),
    @("sample4", 5.5, 5, "comment 4")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$ws.Range("C3").Select()

$wb.Windows.Item(1).Left = 2790
